$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E17").Value = 90
$ws.Range("E18").Value = 135
$ws.Range("E20").Value = 180
$ws.Range("E21").Value = 90
$ws.Range("E36").Value = 180
